$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value2 = 468.85715
$ws.Range("I41").Value2 = 400
$ws.Range("J41").Value2 = 496.4
$ws.Range("K41").Value2 = 400
$ws.Range("L41").Value2 = 496.4
$ws.Range("M41").Value2 = 40
$ws.Range("N41").Value2 = -1376.4
$ws.Range("H62").Value2 = 2244.4443
$ws.Range("I62").Value2 = 1433.3334
$ws.Range("J62").Value2 = 2650
$ws.Range("K62").Value2 = 1433.3334
$ws.Range("L62").Value2 = 2650
$ws.Range("M62").Value2 = -809.3334
$ws.Range("N62").Value2 = -3898
$ws.Range("H65").Value2 = 2244.4443
$ws.Range("I65").Value2 = 1433.3334
$ws.Range("J65").Value2 = 2650
$ws.Range("K65").Value2 = 7166.666999999999
$ws.Range("L65").Value2 = 13250
$ws.Range("M65").Value2 = -4046.666999999999
$ws.Range("N65").Value2 = -19490
$ws.Range("H86").Value2 = 3000
$ws.Range("I86").Value2 = 3000
$ws.Range("J86").Value2 = 3000
$ws.Range("K86").Value2 = 3000
$ws.Range("L86").Value2 = 3000
$ws.Range("M86").Value2 = -1877
$ws.Range("N86").Value2 = -5246
$ws.Range("H87").Value2 = 19913.334
$ws.Range("J87").Value2 = 19913.334
$ws.Range("L87").Value2 = 19913.334
$ws.Range("N87").Value2 = -22409.334
$ws.Range("H89").Value2 = 3000
$ws.Range("I89").Value2 = 3000
$ws.Range("J89").Value2 = 3000
$ws.Range("K89").Value2 = 15000
$ws.Range("L89").Value2 = 15000
$ws.Range("M89").Value2 = -9384
$ws.Range("N89").Value2 = -26232
$ws.Range("H90").Value2 = 19913.334
$ws.Range("J90").Value2 = 19913.334
$ws.Range("L90").Value2 = 59740.00199999999
$ws.Range("N90").Value2 = -72220.00199999999
$ws.Range("H125").Value2 = 14407.5
$ws.Range("I125").Value2 = 21552
$ws.Range("J125").Value2 = 2500
$ws.Range("K125").Value2 = 193968
$ws.Range("L125").Value2 = 22500
$ws.Range("M125").Value2 = -191508
$ws.Range("N125").Value2 = -27420
$ws.Range("H132").Value2 = 4314714
$ws.Range("I132").Value2 = 4468803.5
$ws.Range("J132").Value2 = 206
$ws.Range("K132").Value2 = 13406410.5
$ws.Range("L132").Value2 = 618
$ws.Range("M132").Value2 = -13403880.5
$ws.Range("N132").Value2 = -5678
$ws.Range("H137").Value2 = 2808.5833
$ws.Range("I137").Value2 = 2920.4
$ws.Range("J137").Value2 = 2249.5
$ws.Range("K137").Value2 = 8761.200000000001
$ws.Range("L137").Value2 = 6748.5
$ws.Range("M137").Value2 = -6211.200000000001
$ws.Range("N137").Value2 = -11848.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value2 = 0
$ws.Range("I13").Value2 = 0
$ws.Range("K13").Value2 = 0
$ws.Range("M13").ClearContents()
$ws.Range("H32").Value2 = 7502.95
$ws.Range("I32").Value2 = 5140.1665
$ws.Range("K32").Value2 = 5140.1665
$ws.Range("M32").Value2 = -4853.1665
$ws.Range("H45").Value2 = 1528.2609
$ws.Range("I45").Value2 = 1621.5714
$ws.Range("K45").Value2 = 1621.5714
$ws.Range("M45").Value2 = -1244.5714
$ws.Range("H110").Value2 = 631
$ws.Range("I110").Value2 = 648
$ws.Range("J110").Value2 = 580
$ws.Range("K110").Value2 = 648
$ws.Range("L110").Value2 = 580
$ws.Range("M110").Value2 = 1397
$ws.Range("N110").Value2 = -4670
$ws.Range("H132").Value2 = 399718.7
$ws.Range("I132").Value2 = 48891.41
$ws.Range("J132").Value2 = 1686085.4
$ws.Range("K132").Value2 = 146674.23
$ws.Range("L132").Value2 = 5058256.199999999
$ws.Range("M132").Value2 = -144144.23
$ws.Range("N132").Value2 = -5063316.199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value2 = 9375.429
$ws.Range("I25").Value2 = 432
$ws.Range("J25").Value2 = 21300
$ws.Range("K25").Value2 = 432
$ws.Range("L25").Value2 = 21300
$ws.Range("M25").Value2 = -197
$ws.Range("N25").Value2 = -21770
$ws.Range("H80").Value2 = 1900
$ws.Range("I80").Value2 = 1130
$ws.Range("J80").Value2 = 2175
$ws.Range("K80").Value2 = 1130
$ws.Range("L80").Value2 = 2175
$ws.Range("M80").Value2 = -132
$ws.Range("N80").Value2 = -4171
$ws.Range("H83").Value2 = 1900
$ws.Range("I83").Value2 = 1130
$ws.Range("J83").Value2 = 2175
$ws.Range("K83").Value2 = 5650
$ws.Range("L83").Value2 = 10875
$ws.Range("M83").Value2 = -658
$ws.Range("N83").Value2 = -20859
$ws.Range("H134").Value2 = 130329.625
$ws.Range("I134").Value2 = 252405.75
$ws.Range("J134").Value2 = 8253.5
$ws.Range("K134").Value2 = 757217.25
$ws.Range("L134").Value2 = 24760.5
$ws.Range("M134").Value2 = -754682.25
$ws.Range("N134").Value2 = -29830.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1001.6111
$ws.Range("I16").Value2 = 985.75
$ws.Range("J16").Value2 = 1033.3334
$ws.Range("K16").Value2 = 985.75
$ws.Range("L16").Value2 = 1033.3334
$ws.Range("M16").Value2 = -698.75
$ws.Range("N16").Value2 = -1607.3334
$ws.Range("H31").Value2 = 53736.08
$ws.Range("I31").Value2 = 69450.63
$ws.Range("K31").Value2 = 69450.63
$ws.Range("M31").Value2 = -69155.63
$ws.Range("H34").Value2 = 53736.08
$ws.Range("I34").Value2 = 69450.63
$ws.Range("K34").Value2 = 69450.63
$ws.Range("M34").Value2 = -69248.63
$ws.Range("H113").Value2 = 1001.6111
$ws.Range("I113").Value2 = 985.75
$ws.Range("J113").Value2 = 1033.3334
$ws.Range("K113").Value2 = 985.75
$ws.Range("L113").Value2 = 1033.3334
$ws.Range("M113").Value2 = 1184.25
$ws.Range("N113").Value2 = -5373.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value2 = 1610
$ws.Range("I25").Value2 = 1262.5
$ws.Range("J25").Value2 = 3000
$ws.Range("K25").Value2 = 3787.5
$ws.Range("L25").Value2 = 9000
$ws.Range("M25").Value2 = -3618.5
$ws.Range("N25").Value2 = -9338
$ws.Range("H30").Value2 = 1610
$ws.Range("I30").Value2 = 1262.5
$ws.Range("J30").Value2 = 3000
$ws.Range("K30").Value2 = 3787.5
$ws.Range("L30").Value2 = 9000
$ws.Range("M30").Value2 = -3685.5
$ws.Range("N30").Value2 = -9204
$ws.Range("H34").Value2 = 1357.375
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 1357.375
$ws.Range("K34").Value2 = 0
$ws.Range("L34").Value2 = 4072.125
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value2 = -4240.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2640.6
$ws.Range("I80").Value2 = 2550
$ws.Range("J80").Value2 = 2776.5
$ws.Range("K80").Value2 = 2550
$ws.Range("L80").Value2 = 2776.5
$ws.Range("M80").Value2 = -1552
$ws.Range("N80").Value2 = -4772.5
$ws.Range("H83").Value2 = 2640.6
$ws.Range("I83").Value2 = 2550
$ws.Range("J83").Value2 = 2776.5
$ws.Range("K83").Value2 = 12750
$ws.Range("L83").Value2 = 13882.5
$ws.Range("M83").Value2 = -7758
$ws.Range("N83").Value2 = -23866.5
$ws.Range("H107").Value2 = 1088.0769
$ws.Range("I107").Value2 = 850
$ws.Range("J107").Value2 = 1159.5
$ws.Range("K107").Value2 = 850
$ws.Range("L107").Value2 = 1159.5
$ws.Range("M107").Value2 = 1070
$ws.Range("N107").Value2 = -4999.5
$ws.Range("H132").Value2 = 48697.625
$ws.Range("I132").Value2 = 73950
$ws.Range("J132").Value2 = 13344.3
$ws.Range("K132").Value2 = 221850
$ws.Range("L132").Value2 = 40032.89999999999
$ws.Range("M132").Value2 = -219320
$ws.Range("N132").Value2 = -45092.89999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 2172.5293
$ws.Range("I61").Value2 = 2173.111
$ws.Range("J61").Value2 = 2171.875
$ws.Range("K61").Value2 = 2173.111
$ws.Range("L61").Value2 = 2171.875
$ws.Range("M61").Value2 = -1971.111
$ws.Range("N61").Value2 = -2575.875
$ws.Range("H96").Value2 = 30000
$ws.Range("J96").Value2 = 30000
$ws.Range("L96").Value2 = 30000
$ws.Range("N96").Value2 = -35492
$ws.Range("H113").Value2 = 2172.5293
$ws.Range("I113").Value2 = 2173.111
$ws.Range("J113").Value2 = 2171.875
$ws.Range("K113").Value2 = 2173.111
$ws.Range("L113").Value2 = 2171.875
$ws.Range("M113").Value2 = -3.110999999999876
$ws.Range("N113").Value2 = -6511.875
$ws.Range("H132").Value2 = 9419.5
$ws.Range("I132").Value2 = 11820.4
$ws.Range("J132").Value2 = 5418
$ws.Range("K132").Value2 = 35461.2
$ws.Range("L132").Value2 = 16254
$ws.Range("M132").Value2 = -32931.2
$ws.Range("N132").Value2 = -21314

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value2 = 0
$ws.Range("J86").Value2 = 0
$ws.Range("L86").Value2 = 0
$ws.Range("N86").ClearContents()
$ws.Range("H87").Value2 = 24225
$ws.Range("J87").Value2 = 24225
$ws.Range("L87").Value2 = 24225
$ws.Range("N87").Value2 = -26721
$ws.Range("H89").Value2 = 0
$ws.Range("J89").Value2 = 0
$ws.Range("L89").Value2 = 0
$ws.Range("N89").ClearContents()
$ws.Range("H90").Value2 = 24225
$ws.Range("J90").Value2 = 24225
$ws.Range("L90").Value2 = 72675
$ws.Range("N90").Value2 = -85155
$ws.Range("H99").Value2 = 42264.445
$ws.Range("J99").Value2 = 42264.445
$ws.Range("L99").Value2 = 42264.445
$ws.Range("N99").Value2 = -48254.445
$ws.Range("H109").Value2 = 35459
$ws.Range("J109").Value2 = 35459
$ws.Range("L109").Value2 = 35459
$ws.Range("N109").Value2 = -38233
$ws.Range("H132").Value2 = 6988.72
$ws.Range("I132").Value2 = 8042.222
$ws.Range("J132").Value2 = 4279.7144
$ws.Range("K132").Value2 = 24126.666
$ws.Range("L132").Value2 = 12839.1432
$ws.Range("M132").Value2 = -21596.666
$ws.Range("N132").Value2 = -17899.1432
